$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 51 (S_1011, row for social data without a practice run) was missing the
# three "amount/probability" gap columns (J:X and AB:AI) that a practice-run
# record would normally populate. Remove those cells; deleting them with an
# up-shift also realigns the trailing placeholder rows below (they lose their
# now out-of-range Y:AA "NULL" filler cells), matching rows 65-68 collapsing
# into 65-67.
$ws.Range("J51:X51").Delete() | Out-Null

# Columns Y:AA on row 51 remain as blank, styled placeholder cells (like every
# other "no practice run" row), so just clear their contents/type.
$ws.Range("Y51:AA51").ClearContents() | Out-Null

# Reflect the author's new cursor position/view after the cleanup.
$ws.Range("Q52").Select() | Out-Null
